$d = $word.ActiveDocument
Write-Host ("Before count: " + $d.Paragraphs.Count)
$r = $d.Paragraphs.Item(3).Range
$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$xmlFrag = "<w:p $ns><w:pPr><w:pStyle w:val='APATNivel2'/></w:pPr><w:r><w:t>Crear una rama</w:t></w:r><w:r><w:t xml:space='preserve'> (</w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:t>branch</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t>)</w:t></w:r></w:p>"
$r.InsertXML($xmlFrag)
Write-Host ("After count: " + $d.Paragraphs.Count)
Write-Host ("Para 3 text: " + $d.Paragraphs.Item(3).Range.Text)
Write-Host ("Para 4 text: " + $d.Paragraphs.Item(4).Range.Text.Substring(0,40))
